$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values are preserved as exact text (matching original text formatting)
$textCells = @('D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D12', 'D13', 'D15', 'D16', 'D18', 'D19', 'D20', 'D22', 'D24', 'D25', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D44', 'D45', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $textCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range('D2').Value = '29.441.41'
$ws.Range('E2').Value = '  -2.21%  '
$ws.Range('D3').Value = '1.987.43'
$ws.Range('E3').Value = '  -6.16%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '329.75'
$ws.Range('E5').Value = '  -4.83%  '
$ws.Range('D6').Value = '1.006'
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('D7').Value = '0.4969'
$ws.Range('E7').Value = '  -4.42%  '
$ws.Range('D8').Value = '0.4204'
$ws.Range('E8').Value = '  -5.96%  '
$ws.Range('D9').Value = '51.98'
$ws.Range('E9').Value = '  -4.03%  '
$ws.Range('D10').Value = '0.08879'
$ws.Range('E10').Value = '  -5.21%  '
$ws.Range('E11').Value = '  -5.61%  '
$ws.Range('D12').Value = '23.32'
$ws.Range('E12').Value = '  -7.95%  '
$ws.Range('D13').Value = '8.044'
$ws.Range('E13').Value = '  -7.16%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.959.85'
$ws.Range('E14').Value = '  -6.06%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = '6.494'
$ws.Range('E15').Value = '  -6.82%  '
$ws.Range('D16').Value = '96.00'
$ws.Range('E16').Value = '  -6.38%  '
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range('D18').Value = '0.00001104'
$ws.Range('E18').Value = '  -5.69%  '
$ws.Range('D19').Value = '0.06622'
$ws.Range('D20').Value = '19.69'
$ws.Range('E20').Value = '  -8.77%  '
$ws.Range('E21').Value = '  -0.16%  '
$ws.Range('D22').Value = '5.949'
$ws.Range('E22').Value = '  -5.77%  '
$ws.Range('D23').Value = '29.439.62'
$ws.Range('E23').Value = '  -2.32%  '
$ws.Range('D24').Value = '11.84'
$ws.Range('E24').Value = '  -7.25%  '
$ws.Range('D25').Value = '2.280'
$ws.Range('E25').Value = '  -2.43%  '
$ws.Range('D26').Value = '157.26'
$ws.Range('E26').Value = '  -3.63%  '
$ws.Range('D27').Value = '20.54'
$ws.Range('E27').Value = '  -7.36%  '
$ws.Range('D28').Value = '6.518'
$ws.Range('E28').Value = '  -4.07%  '
$ws.Range('D29').Value = '2.333'
$ws.Range('E29').Value = '  -8.57%  '
$ws.Range('D30').Value = '127.76'
$ws.Range('E30').Value = '  -4.75%  '
$ws.Range('D31').Value = '1.050'
$ws.Range('E31').Value = '  -9.12%  '
$ws.Range('D32').Value = '0.09926'
$ws.Range('E32').Value = '  -6.24%  '
$ws.Range('D33').Value = '1.566'
$ws.Range('E33').Value = '  -12.49%  '
$ws.Range('D34').Value = '5.832'
$ws.Range('E34').Value = '  -7.28%  '
$ws.Range('D35').Value = '3.787'
$ws.Range('E35').Value = '  -4.51%  '
$ws.Range('D36').Value = '9.560'
$ws.Range('E36').Value = '  -11.31%  '
$ws.Range('D37').Value = '0.02447'
$ws.Range('E37').Value = '  -7.53%  '
$ws.Range('D38').Value = '0.06345'
$ws.Range('E38').Value = '  -7.76%  '
$ws.Range('D39').Value = '1.284'
$ws.Range('E39').Value = '  -3.76%  '
$ws.Range('B40').Value = 'Aptos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D40').Value = '11.75'
$ws.Range('E40').Value = '  -7.79%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = '0.6506'
$ws.Range('E41').Value = '  -8.79%  '
$ws.Range('E42').Value = '  -8.24%  '
$ws.Range('E43').Value = '  -0.02%  '
$ws.Range('D44').Value = '0.6341'
$ws.Range('E44').Value = '  -9.26%  '
$ws.Range('D45').Value = '2.215'
$ws.Range('E45').Value = '  -7.89%  '
$ws.Range('E46').Value = '  -9.27%  '
$ws.Range('D47').Value = '1.266'
$ws.Range('E47').Value = '  +0.18%  '
$ws.Range('D48').Value = '3.533'
$ws.Range('E48').Value = '  -2.79%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.00000000332'
$ws.Range('E49').Value = '  -5.42%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.06985'
$ws.Range('E50').Value = '  -3.08%  '
$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D51').Value = '1.138'
$ws.Range('E51').Value = '  -6.24%  '
